# Update "想去人数" (interest count) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12896
$ws1.Range("F3").Value  = 7221
$ws1.Range("F12").Value = 361
$ws1.Range("F19").Value = 375
$ws1.Range("F24").Value = 202
$ws1.Range("F26").Value = 5266
$ws1.Range("F28").Value = 1449
$ws1.Range("F29").Value = 317
$ws1.Range("F30").Value = 1722

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 9296

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9296
$ws4.Range("F5").Value  = 12896
$ws4.Range("F6").Value  = 7221
$ws4.Range("F12").Value = 361
$ws4.Range("F18").Value = 375
$ws4.Range("F26").Value = 202
$ws4.Range("F28").Value = 5266
$ws4.Range("F30").Value = 1449
$ws4.Range("F33").Value = 317
$ws4.Range("F35").Value = 1722
